$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 5.25
$ws.Range("H2").Value = 3.3
$ws.Range("I2").Value = 1.75
$ws.Range("L2").Value = 2.5
$ws.Range("Q2").Value = 2.3
$ws.Range("R2").Value = 1.6
$ws.Range("W2").Value = 11
$ws.Range("X2").Value = 23
$ws.Range("Y2").Value = 17
$ws.Range("AI2").Value = 7.5
$ws.Range("AW2").Value = 3.6

# Row 5 updates
$ws.Range("M5").Value = 1.08
$ws.Range("N5").Value = 8

# Row 7 updates
$ws.Range("N7").Value = 10
